$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.032773962162504
$ws.Range("D2").Value = 1.04252438505094
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.051036308637997
$ws.Range("I2").Value = 1.039173507926888
$ws.Range("J2").Value = 1.037902336002627
$ws.Range("K2").Value = 1.045300828908993
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.053788925128773
$ws.Range("N2").Value = 1.016624743186079
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.033506138069271
$ws.Range("D3").Value = 1.04310179610213
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.051744122570948
$ws.Range("I3").Value = 1.039323567396665
$ws.Range("J3").Value = 1.038278051018933
$ws.Range("K3").Value = 1.045689700825592
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.054309574758401
$ws.Range("N3").Value = 1.016749520146816
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033980567715573
$ws.Range("D4").Value = 1.043476017254895
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.052203080366978
$ws.Range("I4").Value = 1.039419842561915
$ws.Range("J4").Value = 1.038521120052491
$ws.Range("K4").Value = 1.045941221735663
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.054646748097956
$ws.Range("N4").Value = 1.016830230077793
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.034180174619826
$ws.Range("D5").Value = 1.043633480883893
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.052396252573717
$ws.Range("I5").Value = 1.039460118741043
$ws.Range("J5").Value = 1.038623294481308
$ws.Range("K5").Value = 1.046046934484557
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.054788560581114
$ws.Range("N5").Value = 1.016864153099685
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.034213698626695
$ws.Range("D6").Value = 1.043659927933213
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.052428700245145
$ws.Range("I6").Value = 1.039466869660832
$ws.Range("J6").Value = 1.03864044928727
$ws.Range("K6").Value = 1.046064682515753
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.054812375269401
$ws.Range("N6").Value = 1.016869848474916
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033983234259671
$ws.Range("D7").Value = 1.043478120739887
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.052205660656885
$ws.Range("I7").Value = 1.039420381512789
$ws.Range("J7").Value = 1.038522485360216
$ws.Range("K7").Value = 1.045942634381094
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.054648642750888
$ws.Range("N7").Value = 1.016830683388664
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033021265975224
$ws.Range("D8").Value = 1.042719398877259
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.051275318614051
$ws.Range("I8").Value = 1.039224391231419
$ws.Range("J8").Value = 1.038029318901798
$ws.Range("K8").Value = 1.045432270909014
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.053964822424722
$ws.Range("N8").Value = 1.016666917891414
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.031331320873087
$ws.Range("D9").Value = 1.041387099613741
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.049643349303587
$ws.Range("I9").Value = 1.038872760451234
$ws.Range("J9").Value = 1.037160020665314
$ws.Range("K9").Value = 1.044532203546591
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.052762053750051
$ws.Range("N9").Value = 1.016378139794907
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.030208284571882
$ws.Range("D10").Value = 1.040502157149669
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.048560487666905
$ws.Range("I10").Value = 1.038634171118656
$ws.Range("J10").Value = 1.036580382553114
$ws.Range("K10").Value = 1.043931746792809
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.051961801292344
$ws.Range("N10").Value = 1.016185512805927
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.029722875388648
$ws.Range("D11").Value = 1.040119765338063
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.048092839053954
$ws.Range("I11").Value = 1.038529880556164
$ws.Range("J11").Value = 1.036329384592097
$ws.Range("K11").Value = 1.04367166312877
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.051615683849776
$ws.Range("N11").Value = 1.016102083409967
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.029542706031156
$ws.Range("D12").Value = 1.039977849383968
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.047919321808025
$ws.Range("I12").Value = 1.038490995996362
$ws.Range("J12").Value = 1.036236152514757
$ws.Range("K12").Value = 1.043575045662238
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.051487181944065
$ws.Range("N12").Value = 1.016071091395045
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.029581346907853
$ws.Range("D13").Value = 1.040008285311272
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.047956533289232
$ws.Range("I13").Value = 1.038499343478496
$ws.Range("J13").Value = 1.03625615110159
$ws.Range("K13").Value = 1.043595770902133
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.051514743225225
$ws.Range("N13").Value = 1.016077739398459
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029707979803579
$ws.Range("D14").Value = 1.040108032032659
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.048078492210133
$ws.Range("I14").Value = 1.038526669328214
$ws.Range("J14").Value = 1.036321677996873
$ws.Range("K14").Value = 1.043663676910799
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.051605060576385
$ws.Range("N14").Value = 1.016099521651083
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.029786020225876
$ws.Range("D15").Value = 1.040169505427461
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.048153660143874
$ws.Range("I15").Value = 1.038543486310543
$ws.Range("J15").Value = 1.036362051286319
$ws.Range("K15").Value = 1.043705514655716
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.05166071628518
$ws.Range("N15").Value = 1.016112942082601
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.030240518149868
$ws.Range("D16").Value = 1.040527552143911
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.048591550239625
$ws.Range("I16").Value = 1.038641071964979
$ws.Range("J16").Value = 1.03659704034976
$ws.Range("K16").Value = 1.043949006090583
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.051984780493104
$ws.Range("N16").Value = 1.01619104934617
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030525847867635
$ws.Range("D17").Value = 1.040752359469397
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.048866560121927
$ws.Range("I17").Value = 1.038702023112624
$ws.Range("J17").Value = 1.036744440847151
$ws.Range("K17").Value = 1.044101720956525
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.05218816504887
$ws.Range("N17").Value = 1.016240038781291
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030692359952043
$ws.Range("D18").Value = 1.040883562260003
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.049027087989808
$ws.Range("I18").Value = 1.038737480294723
$ws.Range("J18").Value = 1.036830415926489
$ws.Range("K18").Value = 1.044190788949138
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.052306834014407
$ws.Range("N18").Value = 1.016268611466616
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030749150469308
$ws.Range("D19").Value = 1.040928311930067
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.049081844001271
$ws.Range("I19").Value = 1.038749554208795
$ws.Range("J19").Value = 1.036859730972237
$ws.Range("K19").Value = 1.044221157429601
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.052347303487041
$ws.Range("N19").Value = 1.016278353656174
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030495225976734
$ws.Range("D20").Value = 1.040728231845105
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.048837041812167
$ws.Range("I20").Value = 1.038695493409532
$ws.Range("J20").Value = 1.036728626276502
$ws.Range("K20").Value = 1.0440853369139
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.052166339846767
$ws.Range("N20").Value = 1.016234782884371
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.02967068587767
$ws.Range("D21").Value = 1.040078655741642
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.048042573131763
$ws.Range("I21").Value = 1.038518626581487
$ws.Range("J21").Value = 1.036302381955591
$ws.Range("K21").Value = 1.043643680570267
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.051578462654335
$ws.Range("N21").Value = 1.016093107389153
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02915303662362
$ws.Range("D22").Value = 1.039670944821178
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.047544149569335
$ws.Range("I22").Value = 1.038406576752205
$ws.Range("J22").Value = 1.03603438441172
$ws.Range("K22").Value = 1.043365931872961
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.051209197802289
$ws.Range("N22").Value = 1.016004015444941
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.029427378433525
$ws.Range("D23").Value = 1.039887012729512
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.047808269192437
$ws.Range("I23").Value = 1.038466056469375
$ws.Range("J23").Value = 1.036176454601674
$ws.Range("K23").Value = 1.043513177074102
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.051404917640856
$ws.Range("N23").Value = 1.016051246026679
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.030509062431592
$ws.Range("D24").Value = 1.040739133843934
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.048850379497586
$ws.Range("I24").Value = 1.038698444193671
$ws.Range("J24").Value = 1.036735772203947
$ws.Range("K24").Value = 1.044092740182457
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.05217620159838
$ws.Range("N24").Value = 1.016237157804268
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03176758758705
$ws.Range("D25").Value = 1.041730965614108
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.050064360245509
$ws.Range("I25").Value = 1.038964403689816
$ws.Range("J25").Value = 1.037384779592364
$ws.Range("K25").Value = 1.044764971477386
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.053072725301116
$ws.Range("N25").Value = 1.016452816968958
